$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 300
$ws.Range("I29").Value = 300
$ws.Range("K29").Value = 900
$ws.Range("M29").Value = -619

# Row 33
$ws.Range("H33").Value = 25000650
$ws.Range("I33").Value = 25000650
$ws.Range("K33").Value = 25000650
$ws.Range("M33").Value = -25000421

# Row 38
$ws.Range("H38").Value = 70
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 162
$ws.Range("N38").ClearContents()

# Row 43
$ws.Range("H43").Value = 15757.857
$ws.Range("I43").Value = 50555.5
$ws.Range("J43").Value = 1838.8
$ws.Range("K43").Value = 50555.5
$ws.Range("L43").Value = 1838.8
$ws.Range("M43").Value = -50486.5
$ws.Range("N43").Value = -1976.8

# Row 46
$ws.Range("H46").Value = 6500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 19500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -19738

# Row 58
$ws.Range("H58").Value = 463.16666
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 6000
$ws.Range("N58").Value = -6300

# Row 60
$ws.Range("H60").Value = 6500
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 6500
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 19500
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -20468

# Row 87
$ws.Range("H87").Value = 80250
$ws.Range("J87").Value = 81333.336
$ws.Range("L87").Value = 81333.336
$ws.Range("N87").Value = -83829.336

# Row 90
$ws.Range("H90").Value = 80250
$ws.Range("J90").Value = 81333.336
$ws.Range("L90").Value = 244000.008
$ws.Range("N90").Value = -256480.008

# Row 107
$ws.Range("H107").Value = 1431.6666
$ws.Range("I107").Value = 1318
$ws.Range("K107").Value = 1318
$ws.Range("M107").Value = 602

# Row 112
$ws.Range("H112").Value = 54365.05
$ws.Range("J112").Value = 68622.47
$ws.Range("L112").Value = 205867.41
$ws.Range("N112").Value = -208083.41

# Row 137
$ws.Range("H137").Value = 7574.4
$ws.Range("I137").Value = 1718.125
$ws.Range("J137").Value = 30999.5
$ws.Range("K137").Value = 5154.375
$ws.Range("L137").Value = 92998.5
$ws.Range("M137").Value = -2604.375
$ws.Range("N137").Value = -98098.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 80434.71000000001
$ws.Range("I45").Value = 93389.336
$ws.Range("K45").Value = 93389.336
$ws.Range("M45").Value = -93012.336

# Row 132
$ws.Range("H132").Value = 3678.862
$ws.Range("I132").Value = 3640.9048
$ws.Range("J132").Value = 3778.5
$ws.Range("K132").Value = 10922.7144
$ws.Range("L132").Value = 11335.5
$ws.Range("M132").Value = -8392.714399999999
$ws.Range("N132").Value = -16395.5

$ws = $wb.Worksheets.Item("BSM")
# Row 54
$ws.Range("H54").Value = 36666.332
$ws.Range("I54").Value = 36666.332
$ws.Range("K54").Value = 36666.332
$ws.Range("M54").Value = -36182.332

# Row 107
$ws.Range("H107").Value = 2584.647
$ws.Range("I107").Value = 2366
$ws.Range("J107").Value = 3428
$ws.Range("K107").Value = 2366
$ws.Range("L107").Value = 3428
$ws.Range("M107").Value = -446
$ws.Range("N107").Value = -7268

# Row 134
$ws.Range("H134").Value = 2229.5312
$ws.Range("I134").Value = 2204.6775
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6614.032499999999
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4079.032499999999
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 30837.973
$ws.Range("J31").Value = 3241.6428
$ws.Range("L31").Value = 3241.6428
$ws.Range("N31").Value = -3831.6428

# Row 34
$ws.Range("H34").Value = 30837.973
$ws.Range("J34").Value = 3241.6428
$ws.Range("L34").Value = 3241.6428
$ws.Range("N34").Value = -3645.6428

# Row 58
$ws.Range("H58").Value = 3002.7827
$ws.Range("I58").Value = 2798.0278
$ws.Range("J58").Value = 3739.9
$ws.Range("K58").Value = 2798.0278
$ws.Range("L58").Value = 3739.9
$ws.Range("M58").Value = -2595.0278
$ws.Range("N58").Value = -4145.9

# Row 111
$ws.Range("H111").Value = 13450.333
$ws.Range("J111").Value = 20702
$ws.Range("L111").Value = 20702
$ws.Range("N111").Value = -28882

# Row 136
$ws.Range("H136").Value = 3002.7827
$ws.Range("I136").Value = 2798.0278
$ws.Range("J136").Value = 3739.9
$ws.Range("K136").Value = 8394.0834
$ws.Range("L136").Value = 11219.7
$ws.Range("M136").Value = -5844.0834
$ws.Range("N136").Value = -16319.7

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 42540430
$ws.Range("I4").Value = 45735852
$ws.Range("J4").Value = 999999.5
$ws.Range("K4").Value = 137207556
$ws.Range("L4").Value = 2999998.5
$ws.Range("M4").Value = -137207444
$ws.Range("N4").Value = -3000222.5

# Row 5
$ws.Range("H5").Value = 807.36365
$ws.Range("J5").Value = 1149.75
$ws.Range("L5").Value = 3449.25
$ws.Range("N5").Value = -3673.25

# Row 23
$ws.Range("H23").Value = 741.75
$ws.Range("I23").Value = 498
$ws.Range("K23").Value = 1494
$ws.Range("M23").Value = -1259

# Row 24
$ws.Range("H24").Value = 2500
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2500
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 7500
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -7960

# Row 104
$ws.Range("H104").Value = 8298
$ws.Range("J104").Value = 10847
$ws.Range("L104").Value = 32541
$ws.Range("N104").Value = -37783

# Row 135
$ws.Range("H135").Value = 807.36365
$ws.Range("J135").Value = 1149.75
$ws.Range("L135").Value = 10347.75
$ws.Range("N135").Value = -15417.75

# Row 139
$ws.Range("H139").Value = 1480.1904
$ws.Range("I139").Value = 1226.8889
$ws.Range("K139").Value = 3680.6667
$ws.Range("M139").Value = 1459.3333

$ws = $wb.Worksheets.Item("GSM")
# Row 23
$ws.Range("H23").Value = 1723.1428
$ws.Range("I23").Value = 906
$ws.Range("J23").Value = 2050
$ws.Range("K23").Value = 906
$ws.Range("L23").Value = 2050
$ws.Range("M23").Value = -683
$ws.Range("N23").Value = -2496

# Row 29
$ws.Range("H29").Value = 258248.75
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 70
$ws.Range("H70").Value = 25819.666
$ws.Range("J70").Value = 25819.666
$ws.Range("L70").Value = 25819.666
$ws.Range("N70").Value = -26359.666

# Row 73
$ws.Range("H73").Value = 25819.666
$ws.Range("J73").Value = 25819.666
$ws.Range("L73").Value = 25819.666
$ws.Range("N73").Value = -27691.666

# Row 102
$ws.Range("H102").Value = 41667836
$ws.Range("I102").Value = 1219.9565
$ws.Range("K102").Value = 1219.9565
$ws.Range("M102").Value = 402.0435

# Row 126
$ws.Range("H126").Value = 21821.084
$ws.Range("I126").Value = 27483.666
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 82450.99800000001
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -79980.99800000001
$ws.Range("N126").Value = -19440.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 61
$ws.Range("H61").Value = 3000.158
$ws.Range("I61").Value = 2966.8438
$ws.Range("K61").Value = 2966.8438
$ws.Range("M61").Value = -2764.8438

# Row 113
$ws.Range("H113").Value = 3000.158
$ws.Range("I113").Value = 2966.8438
$ws.Range("K113").Value = 2966.8438
$ws.Range("M113").Value = -796.8438000000001

# Row 132
$ws.Range("H132").Value = 5663.8823
$ws.Range("I132").Value = 5481.8335
$ws.Range("J132").Value = 6100.8
$ws.Range("K132").Value = 16445.5005
$ws.Range("L132").Value = 18302.4
$ws.Range("M132").Value = -13915.5005
$ws.Range("N132").Value = -23362.4

# Row 140
$ws.Range("H140").Value = 209333.33
$ws.Range("J140").Value = 209333.33
$ws.Range("L140").Value = 209333.33
$ws.Range("N140").Value = -219693.33

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 6750
$ws.Range("J18").Value = 7100
$ws.Range("L18").Value = 7100
$ws.Range("N18").Value = -7446

# Row 26
$ws.Range("H26").Value = 9000
$ws.Range("J26").Value = 9000
$ws.Range("L26").Value = 9000
$ws.Range("N26").Value = -9586

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Row 99
$ws.Range("H99").Value = 40000
$ws.Range("J99").Value = 40000
$ws.Range("L99").Value = 40000
$ws.Range("N99").Value = -45990

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 136
$ws.Range("H136").Value = 1278.4468
$ws.Range("I136").Value = 1081.1316
$ws.Range("J136").Value = 2111.5557
$ws.Range("K136").Value = 3243.3948
$ws.Range("L136").Value = 6334.6671
$ws.Range("M136").Value = -693.3948
$ws.Range("N136").Value = -11434.6671

Write-Output "Applied 47 row updates across 8 sheets"